{"js": "// Replace the supervisor's job title.\nconst titleResults = context.document.body.search(\n  \"\u0417\u0430\u043c\u0435\u0441\u0442\u0438\u0442\u0435\u043b\u044c \u0434\u0438\u0440\u0435\u043a\u0442\u043e\u0440\u0430 \u043e\u0442\u0434\u0435\u043b\u0430 \u0438\u043d\u0444\u043e\u0440\u043c\u0430\u0446\u0438\u043e\u043d\u043d\u044b\u0445 \u0442\u0435\u0445\u043d\u043e\u043b\u043e\u0433\u0438\u0439 \",\n  { matchCase: true, matchWholeWord: false }\n);\ntitleResults.load(\"items\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(\n    \"\u0414\u0438\u0440\u0435\u043a\u0442\u043e\u0440 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043f\u0435\u0440\u0441\u043e\u043d\u0430\u043b\u043e\u043c \u0438 \u0432\u043d\u0435\u0448\u043d\u0438\u043c \u0441\u0432\u044f\u0437\u044f\u043c \",\n    Word.InsertLocation.replace\n  );\n}\n\n// Replace the supervisor's full name.\nconst nameResults = context.document.body.search(\n  \"\u0414\u0451\u043c\u0438\u043d \u0410\u0440\u0442\u0451\u043c \u0414\u043c\u0438\u0442\u0440\u0438\u0435\u0432\u0438\u0447\",\n  { matchCase: true, matchWholeWord: false }\n);\nnameResults.load(\"items\");\nawait context.sync();\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\n    \"\u041a\u0430\u043b\u0438\u043d\u0438\u0447\u0435\u043d\u043a\u043e \u0418\u0432\u0430\u043d \u041e\u043b\u0435\u0433\u043e\u0432\u0438\u0447\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n\n// Drop the stale \"_GoBack\" bookmark left over from the previous save.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Replace the supervisor's job title.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u0417\u0430\u043c\u0435\u0441\u0442\u0438\u0442\u0435\u043b\u044c \u0434\u0438\u0440\u0435\u043a\u0442\u043e\u0440\u0430 \u043e\u0442\u0434\u0435\u043b\u0430 \u0438\u043d\u0444\u043e\u0440\u043c\u0430\u0446\u0438\u043e\u043d\u043d\u044b\u0445 \u0442\u0435\u0445\u043d\u043e\u043b\u043e\u0433\u0438\u0439 \"\n$find.Replacement.Text = \"\u0414\u0438\u0440\u0435\u043a\u0442\u043e\u0440 \u043f\u043e \u0443\u043f\u0440\u0430\u0432\u043b\u0435\u043d\u0438\u044e \u043f\u0435\u0440\u0441\u043e\u043d\u0430\u043b\u043e\u043c \u0438 \u0432\u043d\u0435\u0448\u043d\u0438\u043c \u0441\u0432\u044f\u0437\u044f\u043c \"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n\n# Replace the supervisor's full name.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"\u0414\u0451\u043c\u0438\u043d \u0410\u0440\u0442\u0451\u043c \u0414\u043c\u0438\u0442\u0440\u0438\u0435\u0432\u0438\u0447\"\n$find2.Replacement.Text = \"\u041a\u0430\u043b\u0438\u043d\u0438\u0447\u0435\u043d\u043a\u043e \u0418\u0432\u0430\u043d \u041e\u043b\u0435\u0433\u043e\u0432\u0438\u0447\"\n$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n\n# Drop the stale \"_GoBack\" bookmark left over from the previous save.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
